# Updated legacy GSC export data:
# The oldest date's row (2025-10-18) is dropped from the HTTPS export,
# so every later row shifts up by one day/position.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Removing the first data row (row 2) shifts all subsequent rows up by
# one, which reproduces the date-list truncation and count shift seen
# in the diff, and shrinks the used range from A1:C91 to A1:C90.
$ws.Rows.Item(2).Delete()
